$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 405.75
$ws.Range("I6").Value = 1023
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 3069
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = -2957
$ws.Range("N6").Value = -824

$ws.Range("H11").Value = 41
$ws.Range("I11").Value = 41
$ws.Range("K11").Value = 41
$ws.Range("M11").Value = 99

$ws.Range("H19").Value = 1593.375
$ws.Range("I19").Value = 1749.5
$ws.Range("J19").Value = 1541.3334
$ws.Range("K19").Value = 1749.5
$ws.Range("L19").Value = 1541.3334
$ws.Range("M19").Value = -1574.5
$ws.Range("N19").Value = -1891.3334

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H33").Value = 270
$ws.Range("I33").Value = 188
$ws.Range("K33").Value = 188
$ws.Range("M33").Value = 41

$ws.Range("H38").Value = 518.375
$ws.Range("I38").Value = 79.40000000000001
$ws.Range("K38").Value = 238.2
$ws.Range("M38").Value = 133.8

$ws.Range("H39").Value = 259.66666
$ws.Range("I39").Value = 7.8
$ws.Range("J39").Value = 574.5
$ws.Range("K39").Value = 23.4
$ws.Range("L39").Value = 1723.5
$ws.Range("M39").Value = 272.6
$ws.Range("N39").Value = -2315.5

$ws.Range("H62").Value = 4833
$ws.Range("I62").Value = 4833
$ws.Range("K62").Value = 4833
$ws.Range("M62").Value = -4209

$ws.Range("H65").Value = 4833
$ws.Range("I65").Value = 4833
$ws.Range("K65").Value = 24165
$ws.Range("M65").Value = -21045

$ws.Range("H98").Value = 2231.75
$ws.Range("I98").Value = 1914.8572
$ws.Range("J98").Value = 4450
$ws.Range("K98").Value = 1914.8572
$ws.Range("L98").Value = 4450
$ws.Range("M98").Value = -416.8571999999999
$ws.Range("N98").Value = -7446

$ws.Range("H113").Value = 2533
$ws.Range("I113").Value = 2533
$ws.Range("K113").Value = 2533
$ws.Range("M113").Value = 721

$ws.Range("H114").Value = 99995
$ws.Range("J114").Value = 99995
$ws.Range("L114").Value = 99995
$ws.Range("N114").Value = -108673

$ws.Range("H122").Value = 2231.75
$ws.Range("I122").Value = 1914.8572
$ws.Range("J122").Value = 4450
$ws.Range("K122").Value = 5744.571599999999
$ws.Range("L122").Value = 13350
$ws.Range("M122").Value = -3294.571599999999
$ws.Range("N122").Value = -18250

$ws.Range("H135").Value = 2649.6667
$ws.Range("I135").Value = 2499.5
$ws.Range("J135").Value = 2950
$ws.Range("K135").Value = 22495.5
$ws.Range("L135").Value = 26550
$ws.Range("M135").Value = -19960.5
$ws.Range("N135").Value = -31620

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2377.6667
$ws.Range("I20").Value = 2556.6365
$ws.Range("J20").Value = 409
$ws.Range("K20").Value = 2556.6365
$ws.Range("L20").Value = 409
$ws.Range("M20").Value = -2309.6365
$ws.Range("N20").Value = -903

$ws.Range("H86").Value = 787.5
$ws.Range("I86").Value = 685.7143
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 685.7143
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = 437.2857
$ws.Range("N86").Value = -3746

$ws.Range("H89").Value = 787.5
$ws.Range("I89").Value = 685.7143
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 3428.5715
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = 2187.4285
$ws.Range("N89").Value = -18732

$ws.Range("H100").Value = 19500
$ws.Range("J100").Value = 19500
$ws.Range("L100").Value = 19500
$ws.Range("N100").Value = -21664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 12541.25
$ws.Range("J41").Value = 13388.333
$ws.Range("L41").Value = 13388.333
$ws.Range("N41").Value = -14244.333

$ws.Range("H59").Value = 29766.666
$ws.Range("J59").Value = 29766.666
$ws.Range("L59").Value = 29766.666
$ws.Range("N59").Value = -32056.666

$ws.Range("H60").Value = 17200
$ws.Range("I60").Value = 14333.333
$ws.Range("K60").Value = 14333.333
$ws.Range("M60").Value = -13822.333

$ws.Range("H74").Value = 21157
$ws.Range("J74").Value = 21157
$ws.Range("L74").Value = 21157
$ws.Range("N74").Value = -22905

$ws.Range("H77").Value = 21157
$ws.Range("J77").Value = 21157
$ws.Range("L77").Value = 63471
$ws.Range("N77").Value = -72207

$ws.Range("H88").Value = 15000
$ws.Range("J88").Value = 15000
$ws.Range("L88").Value = 15000
$ws.Range("N88").Value = -15812

$ws.Range("H91").Value = 15000
$ws.Range("J91").Value = 15000
$ws.Range("L91").Value = 15000
$ws.Range("N91").Value = -17808

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H96").Value = 4027.1
$ws.Range("J96").Value = 4027.1
$ws.Range("L96").Value = 4027.1
$ws.Range("N96").Value = -9519.1

$ws.Range("H122").Value = 1446.5
$ws.Range("I122").Value = 1494
$ws.Range("J122").Value = 1399
$ws.Range("K122").Value = 4482
$ws.Range("L122").Value = 4197
$ws.Range("M122").Value = -2032
$ws.Range("N122").Value = -9097

$ws.Range("H134").Value = 1333
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 833

$ws.Range("H117").Value = 734.8570999999999
$ws.Range("I117").Value = 512.25
$ws.Range("J117").Value = 1031.6666
$ws.Range("K117").Value = 1536.75
$ws.Range("L117").Value = 3094.9998
$ws.Range("M117").Value = 1905.25
$ws.Range("N117").Value = -9978.9998

$ws.Range("H131").Value = 2700
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 2700
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 8100
$ws.Range("N131").Value = -18180
$ws.Range("M131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 13333
$ws.Range("I46").Value = 9999
$ws.Range("J46").Value = 15000
$ws.Range("K46").Value = 9999
$ws.Range("L46").Value = 15000
$ws.Range("N46").Value = -15312
$ws.Range("M46").Value = -9843

$ws.Range("H57").Value = 15513.75
$ws.Range("I57").Value = 2055
$ws.Range("J57").Value = 20000
$ws.Range("K57").Value = 2055
$ws.Range("L57").Value = 20000
$ws.Range("M57").Value = -1235
$ws.Range("N57").Value = -21640

$ws.Range("H80").Value = 3754.8572
$ws.Range("I80").Value = 3371.75
$ws.Range("J80").Value = 4265.6665
$ws.Range("K80").Value = 3371.75
$ws.Range("L80").Value = 4265.6665
$ws.Range("M80").Value = -2373.75
$ws.Range("N80").Value = -6261.6665

$ws.Range("H83").Value = 3754.8572
$ws.Range("I83").Value = 3371.75
$ws.Range("J83").Value = 4265.6665
$ws.Range("K83").Value = 16858.75
$ws.Range("L83").Value = 21328.3325
$ws.Range("M83").Value = -11866.75
$ws.Range("N83").Value = -31312.3325

$ws.Range("H102").Value = 7108.875
$ws.Range("I102").Value = 6839.143
$ws.Range("K102").Value = 6839.143
$ws.Range("M102").Value = -5217.143

$ws.Range("H132").Value = 3276.4
$ws.Range("I132").Value = 3276.4
$ws.Range("K132").Value = 9829.200000000001
$ws.Range("M132").Value = -7299.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1448.4286
$ws.Range("I22").Value = 1099.75
$ws.Range("J22").Value = 1913.3334
$ws.Range("K22").Value = 1099.75
$ws.Range("L22").Value = 1913.3334
$ws.Range("M22").Value = -804.75
$ws.Range("N22").Value = -2503.3334

$ws.Range("H27").Value = 1448.4286
$ws.Range("I27").Value = 1099.75
$ws.Range("J27").Value = 1913.3334
$ws.Range("K27").Value = 1099.75
$ws.Range("L27").Value = 1913.3334
$ws.Range("M27").Value = -992.75
$ws.Range("N27").Value = -2127.3334

$ws.Range("H46").Value = 3246.2
$ws.Range("I46").Value = 1125
$ws.Range("K46").Value = 1125
$ws.Range("M46").Value = -937

$ws.Range("H104").Value = 50000
$ws.Range("J104").Value = 50000
$ws.Range("L104").Value = 50000
$ws.Range("N104").Value = -56988

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 21677.666
$ws.Range("J20").Value = 21677.666
$ws.Range("L20").Value = 21677.666
$ws.Range("N20").Value = -22157.666

$ws.Range("H23").Value = 3653.3333
$ws.Range("I23").Value = 474.5
$ws.Range("J23").Value = 10011
$ws.Range("K23").Value = 474.5
$ws.Range("L23").Value = 10011
$ws.Range("M23").Value = -245.5
$ws.Range("N23").Value = -10469

$ws.Range("H101").Value = 15230
$ws.Range("J101").Value = 15230
$ws.Range("L101").Value = 15230
$ws.Range("N101").Value = -21720

$ws.Range("H104").Value = 17789.666
$ws.Range("J104").Value = 17789.666
$ws.Range("L104").Value = 17789.666
$ws.Range("N104").Value = -24777.666

$ws.Range("H113").Value = 1135.5714
$ws.Range("I113").Value = 465
$ws.Range("K113").Value = 1395
$ws.Range("M113").Value = 775

Write-Host "Applied Kraken_Profits updates across all sheets."
